$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the paragraph that holds the rId10 hyperlink (the last
# hyperlink before the trailing bookmark paragraph) by matching the
# hyperlink's address, rather than hard-coding a paragraph index.
# ------------------------------------------------------------------
$targetUrl = "https://developer.android.com/guide/topics/ui/ui-events?hl=es"
$newUrl = "https://github.com/eddydn/AndroidGridLayout"

$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Hyperlinks.Count -gt 0) {
        $h = $p.Range.Hyperlinks(1)
        if ($h.Address -eq $targetUrl) {
            $targetPara = $p
        }
    }
}

$r = $targetPara.Range

# ------------------------------------------------------------------
# Replace that paragraph with: itself (now carrying a paragraph-mark
# rStyle, as Word does when a new paragraph is split off right after
# a hyperlink run) plus a brand new paragraph containing the new
# AndroidGridLayout hyperlink.
# ------------------------------------------------------------------
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body>' +
  '<w:p><w:pPr><w:rPr><w:rStyle w:val="Hipervnculo"/></w:rPr></w:pPr><w:hyperlink r:id="rIdOldLink" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hipervnculo"/></w:rPr><w:t>' + $targetUrl + '</w:t></w:r></w:hyperlink></w:p>' +
  '<w:p><w:hyperlink r:id="rIdNewLink" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hipervnculo"/></w:rPr><w:t>' + $newUrl + '</w:t></w:r></w:hyperlink></w:p>' +
  '</w:body></w:document></pkg:xmlData></pkg:part>' +
  '<pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml"><pkg:xmlData>' +
  '<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">' +
  '<Relationship Id="rIdOldLink" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="' + $targetUrl + '" TargetMode="External"/>' +
  '<Relationship Id="rIdNewLink" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="' + $newUrl + '" TargetMode="External"/>' +
  '</Relationships></pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)

Write-Output ("ParaCount after insert: " + $d.Paragraphs.Count)
